$d = $word.ActiveDocument

function Escape-Xml($text) {
    $escaped = $text -replace '&', '&amp;'
    $escaped = $escaped -replace '<', '&lt;'
    $escaped = $escaped -replace '>', '&gt;'
    return $escaped
}

function Get-ParaIndexForText($searchText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND:" $searchText
        return -1
    }
    return $rng.Paragraphs.Item(1).Index
}

function Strike-Paragraph($searchText) {
    $idx = Get-ParaIndexForText $searchText
    if ($idx -lt 0) { return }
    $trueP = $d.Paragraphs.Item($idx)
    $trueP.Range.Font.StrikeThrough = 1
}

# Splits a paragraph's content into two runs: a short leading run (kept as-is,
# no strike) and the remaining text run (struck through). Also marks the
# paragraph-mark (pPr/rPr) as struck through, matching Word's own behaviour
# when the whole paragraph is selected and Strikethrough is toggled.
function Strike-Paragraph-KeepLeading($fullText, $leadingLen) {
    $idx = Get-ParaIndexForText $fullText
    if ($idx -lt 0) { return }
    $trueP = $d.Paragraphs.Item($idx)

    # Strike the whole paragraph first (sets pPr/rPr + run rPr uniformly).
    $trueP.Range.Font.StrikeThrough = 1

    $trueP2 = $d.Paragraphs.Item($idx)
    $pStart = $trueP2.Range.Start
    $pEnd = $trueP2.Range.End

    $leadingText = $fullText.Substring(0, $leadingLen)
    $restText = $fullText.Substring($leadingLen)

    $leadingXml = Escape-Xml $leadingText
    $restXml = Escape-Xml $restText

    $leadingSpaceAttr = ""
    if ($leadingText -ne $leadingText.Trim()) { $leadingSpaceAttr = ' xml:space="preserve"' }
    $restSpaceAttr = ""
    if ($restText -ne $restText.Trim()) { $restSpaceAttr = ' xml:space="preserve"' }

    $contentRange = $d.Range($pStart, $pEnd - 1)
    $xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t' + $leadingSpaceAttr + '>' + $leadingXml + '</w:t></w:r><w:r><w:rPr><w:strike/><w:lang w:val="es-AR"/></w:rPr><w:t' + $restSpaceAttr + '>' + $restXml + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $contentRange.InsertXML($xmlFrag)
}

# Handles the "Estos botones..." paragraph: strikes the first sentence,
# splits off the remainder into a brand-new (non-struck) paragraph, and
# rewrites that remainder as "Al " + "presionar ... acceda al Home" (two
# runs), replacing the original "que al " wording with "Al ".
function Split-EstosBotones() {
    $fullText = "Estos botones tienen que completar los campos de email y contraseña con un usuario valido que al presionar  " + [char]0x2018 + "Ingresar" + [char]0x2019 + ", acceda al Home"
    $idx = Get-ParaIndexForText $fullText
    if ($idx -lt 0) { return }
    $trueP = $d.Paragraphs.Item($idx)

    # Strike the whole paragraph first (sets pPr/rPr + run rPr uniformly).
    $trueP.Range.Font.StrikeThrough = 1

    $trueP2 = $d.Paragraphs.Item($idx)
    $pStart = $trueP2.Range.Start
    $pEnd = $trueP2.Range.End

    $splitMarker = "valido "
    $splitPos = $fullText.IndexOf($splitMarker) + $splitMarker.Length
    $firstText = $fullText.Substring(0, $splitPos)                 # "...un usuario valido "
    $secondOld = $fullText.Substring($splitPos)                    # "que al presionar  '...', acceda al Home"

    # Break the split point into its own paragraph.
    $breakPoint = $pStart + $splitPos
    $breakRange = $d.Range($breakPoint, $breakPoint)
    $breakRange.InsertParagraphAfter()

    # Re-fetch paragraphs after the split.
    $firstP = $d.Paragraphs.Item($idx)
    $secondP = $d.Paragraphs.Item($idx + 1)

    Write-Host "firstP text after split:" $firstP.Range.Text
    Write-Host "secondP text after split:" $secondP.Range.Text

    # Remove the leading "que al " from the new second paragraph and replace
    # it with "Al ", rebuilding the run structure as two runs (no strike).
    # Replace the *entire* second paragraph (pPr + runs, including the
    # paragraph mark) in one shot so the inherited strike on the paragraph
    # mark (picked up from InsertParagraphAfter) is not left behind.
    $removePrefix = "que al "
    $afterPrefixText = $secondOld.Substring($removePrefix.Length)  # "presionar  '...', acceda al Home"
    $afterPrefixSpaceAttr = ""
    if ($afterPrefixText -ne $afterPrefixText.Trim()) { $afterPrefixSpaceAttr = ' xml:space="preserve"' }
    $afterPrefixXml = Escape-Xml $afterPrefixText

    $fullRange = $d.Range($secondP.Range.Start, $secondP.Range.End)
    $xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="4"/></w:numPr><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve">Al </w:t></w:r><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t' + $afterPrefixSpaceAttr + '>' + $afterPrefixXml + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $fullRange.InsertXML($xmlFrag)
}

Split-EstosBotones

# ---- Simple strike-only paragraphs ----
Strike-Paragraph "Tiene acceso a los diferentes juegos"
Strike-Paragraph " Tiene que rutear a cada juego"
Strike-Paragraph "Incorporar módulos y loadchiildren"
Strike-Paragraph "Incorporar juegos"
Strike-Paragraph "Ahorcado"
Strike-Paragraph "Mayor o menor"
Strike-Paragraph " Desde un mazo de carta se va a preguntar si la siguiente es mayor o menor"
Strike-Paragraph " El jugador sumara un punto ante cada carta que adivine"

# ---- Paragraphs that keep a leading space un-struck ----
Strike-Paragraph-KeepLeading " No se deben ingresar datos del teclado" 1
Strike-Paragraph-KeepLeading " Utilizar botones para el ingreso de letras" 1

Write-Host "Done"
